# "article 73 is live"
#
# The dashboard's last row (row 7) carries a rotating set of "blog" widget
# cards (one ser:NN id per card). Publishing article 73 bumps the oldest
# card (ser:69) out of the row and shifts the remaining ones down one slot:
#   I7: ser:69 -> ser:71   (was E7's card)
#   E7: ser:71 -> ser:72   (was C7's card)
#   C7: ser:72 -> ser:73   (brand-new card for the freshly published article)
# B7 (the "video" widget) is left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 71"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 72"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 73"

# Editor's cursor ends up on the cell they just typed the new card into.
[void]$ws.Range("I7").Select()
